$wb = $excel.ActiveWorkbook

# --- Sheet "Version": update Start Time / End Time / Run Time ---
$wsVersion = $wb.Worksheets.Item("Version")
$wsVersion.Range("C10").Value = "22:04:37 03-07-2017 Tuesday Pacific Standard Time"
$wsVersion.Range("C11").Value = "22:05:11 03-07-2017 Tuesday Pacific Standard Time"
$wsVersion.Range("C12").Value = "0:0:34 0 days"

# --- Sheet "3_3": update row 3 & 4 values, delete row 5 ---
$ws33 = $wb.Worksheets.Item("3_3")
$ws33.Range("C3").Value = 0.598
$ws33.Range("D3").Value = 0.621
$ws33.Range("E3").Value = 0.574
$ws33.Range("F3").Value = 0.016
$ws33.Range("H3").Value = "0.587009,0.574672,0.604658,0.593782,0.613932,0.620433,0.586459,0.606221,0.607437,0.581628,0.618585,0.577919,0.61226,0.614407,0.575222,0.615147,0.61811,0.59027,0.616262,0.592118,0.617097,0.615324,0.581065,0.584971,0.620623,0.577641,0.592858,0.615147,0.586635,0.583395,0.611058,0.586268,0.578476,0.611622,0.58887,0.613646,0.578381,0.582838,0.616724,0.574393,0.594441,0.619312,0.620161,0.57438,0.583939,0.617743,0.577926,0.602797,0.612831,0.58209,0.617192,0.587301,0.584787,0.575969,0.61599,0.578741,0.581078,0.611336,0.586092,0.585997,0.613653,0.575494,0.589794,0.617743,0.575406,0.594244,0.608463,0.589427,0.615616,0.599733,0.584882,0.578653,0.613381,0.613286,0.5858,0.594624,0.615698,0.586547,0.610786,0.604101,0.586554,0.610881,0.600012,0.581363,0.612817,0.612362,0.583198,0.612722,0.614869,0.610412,0.579407,0.59044,0.582281,0.607804,0.581995,0.613279,0.588307,0.577919,0.615881,0.590446"

$ws33.Range("C4").Value = 9.836
$ws33.Range("D4").Value = 9.853999999999999
$ws33.Range("E4").Value = 9.815
$ws33.Range("H4").Value = "9.824391,9.83664,9.837665,9.843515,9.827733,9.830247,9.845002,9.843046,9.835709,9.815288,9.837665,9.838311,9.84435,9.824479,9.818073,9.845186,9.844534,9.840356,9.8281,9.83274,9.840172,9.843046,9.832095,9.826246,9.83414,9.847319,9.839337,9.832095,9.825227,9.838127,9.845002,9.841103,9.830512,9.83414,9.841565,9.844629,9.834513,9.82346,9.835722,9.848243,9.845648,9.82717,9.828379,9.841375,9.849642,9.83181,9.819656,9.822251,9.852795,9.839982,9.840444,9.819099,9.842856,9.845845,9.840634,9.824391,9.824024,9.838311,9.849351,9.838039,9.821986,9.837475,9.842679,9.847503,9.826436,9.827081,9.83721,9.849908,9.83776,9.823365,9.82253,9.84361,9.85401,9.8256,9.829303,9.839432,9.846578,9.840077,9.830145,9.826062,9.839337,9.843148,9.832475,9.825043,9.83433,9.848617,9.853827,9.830152,9.828854,9.831069,9.842863,9.846762,9.836368,9.824486,9.834969,9.842489,9.836918,9.828385,9.836368,9.846015"

$ws33.Rows.Item(5).Delete()

# --- Sheet "1_8": update row 3 & 4 values, delete row 5 ---
$ws18 = $wb.Worksheets.Item("1_8")
$ws18.Range("C3").Value = 0.001
$ws18.Range("D3").Value = 0.008999999999999999
$ws18.Range("E3").Value = -0.008
$ws18.Range("H3").Value = "0.001427,0.003375,0.00504,0.002253,0.005589,0.003659,-0.003284,0.000576,0.004479,0.001433,0.005034,-0.007175,0.00504,0.006421,0.008906,-0.003876,-0.003882,-0.003315,0.008111,0.004497,0.001445,0.003388,-0.000534,-0.002735,0.001155,-0.005276,5.2e-05,-0.002187,0.005577,0.00504,0.002271,0.001143,0.00282,0.001427,-0.001354,3.3e-05,-0.002446,-0.001071,0.001957,0.002537,0.003129,0.001445,-0.003025,-0.008001,0.001741,0.005022,0.007784,-0.004968,-0.003001,0.00253,0.003369,0.006125,-0.000799,-0.000775,-0.003284,0.001698,0.0006,-0.001354,0.003665,0.000884,-0.004684,-0.001632,0.00364,-0.006047,-0.000263,0.000859,0.001168,-0.004961,0.004214,-0.001071,-0.004696,-0.001373,0.001704,0.006144,0.009208,0.002265,0.00446,-0.004129,0.001698,0.001414,0.000866,0.003357,0.004485,-0.001077,0.001692,0.003351,-0.002742,0.0006,-0.001064,0.008092,0.00451,5.2e-05,-0.000824,0.002271,-0.001619,0.000866,0.006982,0.006674,0.002253,0.00057"

$ws18.Range("C4").Value = 0.001
$ws18.Range("D4").Value = 0.008999999999999999
$ws18.Range("E4").Value = -0.008999999999999999
$ws18.Range("H4").Value = "3.9e-05,0.008099,-0.001342,-0.007434,0.001698,0.002019,0.002802,0.00282,0.001963,0.001396,0.001421,0.003924,0.002555,0.003085,0.000329,0.004473,-0.002483,-0.001052,0.00808,0.0006,0.0006,0.000329,0.002253,0.0006,2.1e-05,0.002814,0.004485,-0.000528,0.00475,-0.003321,0.000853,0.009196,0.000859,0.001938,0.002259,0.001427,-0.000509,0.004189,0.003665,0.000859,0.002518,0.002524,0.00089,3.9e-05,-0.000516,-0.003303,0.001149,-0.008852,-0.001903,0.003659,0.003129,0.001155,-0.003562,0.000576,-0.001656,0.004183,0.002808,0.001458,0.003394,0.001735,0.000317,-0.000824,-0.001323,0.000341,0.001433,-0.000528,0.004479,-0.004659,-0.006349,0.001692,0.004473,0.000594,-0.000522,-0.004943,4.5e-05,-0.002748,0.005583,0.001963,0.002543,-0.000491,0.002796,-0.002723,0.001439,5.2e-05,-0.001323,3.9e-05,0.00282,0.002814,0.006421,-0.000534,0.003647,0.003079,-0.003845,-0.001077,-0.001071,0.003918,0.000847,0.004491,0.005287,0.006409"

$ws18.Rows.Item(5).Delete()
